# Add season-record columns (Wins / Losses / Ties) to the COL_2017 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - match formatting of the other header cells (bold,
# centered, top-aligned, thin border - same look as the existing headers).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1

# Find the last used row on the sheet (data starts at row 2).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 87  # AD
    $ws.Cells.Item($r, 31).Value = 75  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
